# Revert unit test coverage
# Adds a new row 39 (with the same "time" stamp / decoded fields) to each of
# the four worksheets, mirroring the existing row 38 entries.

$wb = $excel.ActiveWorkbook

$rowsData = @{
    "FE_LFT_#1" = @{
        A = 45825.49508101852
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x64"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 356
        I = 15
    }
    "FE_LFT_#2" = @{
        A = 45825.49508101852
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x78"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 376
        I = 14
    }
    "FE_PLT_#1" = @{
        A = 45825.49508101852
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x6A"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 106
        I = 3
    }
    "FE_PLT_#2" = @{
        A = 45825.49508101852
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x6B"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 107
        I = 3
    }
}

$sheetNames = @("FE_LFT_#1", "FE_LFT_#2", "FE_PLT_#1", "FE_PLT_#2")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $row = $rowsData[$name]

    $ws.Range("A39").Value = $row.A
    $ws.Range("A39").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Range("B39").Value = $row.B
    $ws.Range("C39").Value = $row.C
    $ws.Range("D39").Value = $row.D
    $ws.Range("E39").Value = $row.E
    $ws.Range("F39").Value = $row.F
    $ws.Range("G39").Value = $row.G
    $ws.Range("H39").Value = $row.H
    $ws.Range("I39").Value = $row.I
}

Write-Output "Row 39 added to all sheets"
